$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Task 33 row: update text to include new bug-fix description and mark status Pending
$ws.Range("A33").Value = "Task 33:Bug Fix: clicking edit/delete button leads to wrong id"
$ws.Range("B33").Value = "Pending"
$ws.Range("B33").Font.Bold = $true
$ws.Range("B33").Font.Color = 49407

# Task 32 row: update text to include new bug-fix description and mark status Pending
$ws.Range("A32").Value = "Task 32:Bug Fix: clicking project  edit button immediately leads to wrong id"
$ws.Range("B32").Value = "Pending"
$ws.Range("B32").Font.Bold = $true
$ws.Range("B32").Font.Color = 49407

# Update the active selection to reflect where the author left off editing
$ws.Range("C32").Select()
